$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1721611721611722
$ws.Range("C2").Value = 0.608058608058608
$ws.Range("J2").Value = 0.02564102564102564
$ws.Range("P2").Value = 0.1245421245421245
$ws.Range("S2").Value = 0.0695970695970696
$ws.Range("B3").Value = 0.005952380952380952
$ws.Range("C3").Value = 0.02380952380952381
$ws.Range("J3").Value = 0.02976190476190476
$ws.Range("P3").Value = 0.7440476190476191
$ws.Range("S3").Value = 0.1964285714285714
$ws.Range("J4").Value = 0.05
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.2
$ws.Range("B6").Value = 0.07446808510638298
$ws.Range("D6").Value = 0.01063829787234043
$ws.Range("F6").Value = 0.1223404255319149
$ws.Range("J6").Value = 0.2446808510638298
$ws.Range("O6").Value = 0.02127659574468085
$ws.Range("Q6").Value = 0.1276595744680851
$ws.Range("R6").Value = 0.05851063829787234
$ws.Range("S6").Value = 0.3404255319148936
$ws.Range("B7").Value = 0.1165644171779141
$ws.Range("D7").Value = 0.01226993865030675
$ws.Range("F7").Value = 0.03680981595092025
$ws.Range("J7").Value = 0.1717791411042945
$ws.Range("O7").Value = 0.01840490797546012
$ws.Range("Q7").Value = 0.2331288343558282
$ws.Range("R7").Value = 0.08588957055214724
$ws.Range("S7").Value = 0.3251533742331288
$ws.Range("B8").Value = 0.09826589595375723
$ws.Range("D8").Value = 0.005780346820809248
$ws.Range("F8").Value = 0.06647398843930635
$ws.Range("J8").Value = 0.1011560693641619
$ws.Range("O8").Value = 0.02312138728323699
$ws.Range("Q8").Value = 0.1965317919075145
$ws.Range("R8").Value = 0.1184971098265896
$ws.Range("S8").Value = 0.3901734104046243
$ws.Range("B9").Value = 0.1151832460732984
$ws.Range("D9").Value = 0.03141361256544502
$ws.Range("E9").Value = 0.005235602094240838
$ws.Range("F9").Value = 0.04712041884816754
$ws.Range("J9").Value = 0.162303664921466
$ws.Range("O9").Value = 0.01047120418848168
$ws.Range("Q9").Value = 0.1413612565445026
$ws.Range("R9").Value = 0.05235602094240838
$ws.Range("S9").Value = 0.4345549738219895
$ws.Range("B10").Value = 0.1125211505922166
$ws.Range("D10").Value = 0.02538071065989848
$ws.Range("E10").Value = 0.0008460236886632825
$ws.Range("F10").Value = 0.05752961082910321
$ws.Range("J10").Value = 0.1522842639593909
$ws.Range("O10").Value = 0.01607445008460237
$ws.Range("Q10").Value = 0.2089678510998308
$ws.Range("R10").Value = 0.08629441624365482
$ws.Range("S10").Value = 0.3401015228426396
$ws.Range("G11").Value = 0.150197628458498
$ws.Range("J11").Value = 0.09486166007905138
$ws.Range("K11").Value = 0.2134387351778656
$ws.Range("L11").Value = 0.525691699604743
$ws.Range("S11").Value = 0.0158102766798419
$ws.Range("G12").Value = 0.7391304347826086
$ws.Range("J12").Value = 0.1884057971014493
$ws.Range("K12").Value = 0.007246376811594203
$ws.Range("L12").Value = 0.02173913043478261
$ws.Range("S12").Value = 0.04347826086956522
$ws.Range("F13").Value = 0.02631578947368421
$ws.Range("G13").Value = 0.6052631578947368
$ws.Range("J13").Value = 0.3157894736842105
$ws.Range("S13").Value = 0.05263157894736842
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.01463414634146342
$ws.Range("H15").Value = 0.09268292682926829
$ws.Range("I15").Value = 0.06341463414634146
$ws.Range("J15").Value = 0.3951219512195122
$ws.Range("K15").Value = 0.08292682926829269
$ws.Range("M15").Value = 0.01463414634146342
$ws.Range("O15").Value = 0.1073170731707317
$ws.Range("S15").Value = 0.2292682926829268
$ws.Range("F16").Value = 0.02222222222222222
$ws.Range("H16").Value = 0.1722222222222222
$ws.Range("I16").Value = 0.06666666666666667
$ws.Range("J16").Value = 0.4166666666666667
$ws.Range("K16").Value = 0.1277777777777778
$ws.Range("M16").Value = 0.02222222222222222
$ws.Range("N16").Value = 0.005555555555555556
$ws.Range("O16").Value = 0.1
$ws.Range("S16").Value = 0.06666666666666667
$ws.Range("F17").Value = 0.0175
$ws.Range("H17").Value = 0.1675
$ws.Range("I17").Value = 0.1125
$ws.Range("J17").Value = 0.4125
$ws.Range("K17").Value = 0.1075
$ws.Range("M17").Value = 0.015
$ws.Range("O17").Value = 0.06
$ws.Range("S17").Value = 0.1075
$ws.Range("F18").Value = 0.03389830508474576
$ws.Range("H18").Value = 0.1694915254237288
$ws.Range("I18").Value = 0.06779661016949153
$ws.Range("J18").Value = 0.4463276836158192
$ws.Range("K18").Value = 0.05649717514124294
$ws.Range("M18").Value = 0.01694915254237288
$ws.Range("O18").Value = 0.06214689265536723
$ws.Range("S18").Value = 0.1468926553672316
$ws.Range("F19").Value = 0.0168697282099344
$ws.Range("H19").Value = 0.1808809746954077
$ws.Range("I19").Value = 0.1002811621368322
$ws.Range("J19").Value = 0.3823805060918463
$ws.Range("K19").Value = 0.09746954076850985
$ws.Range("M19").Value = 0.02061855670103093
$ws.Range("O19").Value = 0.06654170571696345
$ws.Range("S19").Value = 0.1349578256794752
